# Updates the cached "datetimeFigureOut" date field text (slide master +
# every slide layout) from the old auto-date to the new one, and relabels
# the sub-figure caption on slide 2 from "(d)" to "(h)".

$p = $ppt.ActivePresentation

$oldDate = [char]0x05D8 + [char]0x0022 + [char]0x05D5 + "/" + [char]0x05D0 + [char]0x05D1 + "/" + [char]0x05EA + [char]0x05E9 + [char]0x05E4 + [char]0x0022 + [char]0x05D2
$newDate = [char]0x05DB + "'" + "/" + [char]0x05D0 + [char]0x05D1 + "/" + [char]0x05EA + [char]0x05E9 + [char]0x05E4 + [char]0x0022 + [char]0x05D2

# Slide master date placeholder.
$master = $p.SlideMaster
$masterDate = $master.Shapes.Item("Date Placeholder 3")
if ($masterDate.TextFrame.TextRange.Text -eq $oldDate) {
    $masterDate.TextFrame.TextRange.Text = $newDate
}

# Every slide layout's date placeholder.
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $shapes = $layout.Shapes
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide 2: relabel the "(d)" sub-figure caption to "(h)".
$slide2 = $p.Slides.Item(2)
$label = $slide2.Shapes.Item("TextBox 48")
if ($label.TextFrame.TextRange.Text -eq "(d)") {
    $label.TextFrame.TextRange.Text = "(h)"
}
